# trafo_id -> gridnode_id refactor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header in column J (row 1) from "trafo_id" to "gridnode_id"
$ws.Range("J1").Value = "gridnode_id"

# Update the active selection to match the new cursor position
$ws.Range("G8").Select()
